$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.734.40'
$ws.Range("E2").Value = '  -2.51%  '
$ws.Range("D3").Value = '3.659.05'
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '''407.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").Value = '''133.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("D7").Value = '3.657.85'
$ws.Range("E7").Value = '  +3.19%  '
$ws.Range("D8").Value = '''0.622'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D10").Value = '''0.730'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.78%  '
$ws.Range("D11").Value = '''0.164'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.02%  '
$ws.Range("D12").Value = '''0.0000332'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.40%  '
$ws.Range("D13").Value = '''42.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '''9.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").Value = '4.246.97'
$ws.Range("E15").Value = '  +3.30%  '
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").Value = '3.683.25'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '''13.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +11.36%  '
$ws.Range("D19").Value = '''20.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '''1.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").Value = '64.807.89'
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("D22").Value = '''422.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.60%  '
$ws.Range("D23").Value = '''15.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +18.37%  '
$ws.Range("D24").Value = '''86.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.78%  '
$ws.Range("D25").Value = '''3.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.27%  '
$ws.Range("D26").Value = '''35.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.96%  '
$ws.Range("D27").Value = '''3.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '''9.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.53%  '
$ws.Range("D29").Value = '''5.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.12%  '
$ws.Range("D30").Value = '''12.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.84%  '
$ws.Range("E31").Value = '  -1.41%  '
$ws.Range("D32").Value = '''0.118'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '''41.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.27%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '''6.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.55%  '
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("D36").Value = '''55.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '''0.0469'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = '''2.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +31.26%  '
$ws.Range("D40").Value = '''0.141'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.19%  '
$ws.Range("D41").Value = '''0.995'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").Value = '0.0₃0667'
$ws.Range("E42").Value = '  -7.11%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").Value = '''4.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.40%  '
$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").Value = '''3.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.78%  '
$ws.Range("D45").Value = '''26.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +27.35%  '
$ws.Range("D46").Value = '''3.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +23.97%  '
$ws.Range("D47").Value = '''2.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.37%  '
$ws.Range("D48").Value = '''144.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("D49").Value = '''2.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.14%  '
$ws.Range("E50").Value = '  -6.67%  '
$ws.Range("D51").Value = '''0.292'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.06%  '
